$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Activity (B2) and ScheduledResource (D2) and PlannedQty (G2)
$ws.Range("B2").Value = "Digital Print 4x0"
$ws.Range("D2").Value = "Digital Press-Labelfire-340"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "6,825"

# Row 3: PlannedQty (G3)
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "3,385"
